$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "50.915.50"
$ws.Range("E2").Value = "  -1.12%  "

$ws.Range("D3").Value = "2.924.21"
$ws.Range("E3").Value = "  -1.66%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'372.39"
$ws.Range("E5").Value = "  -1.58%  "

$ws.Range("D6").Value = "'99.82"
$ws.Range("E6").Value = "  -4.34%  "

$ws.Range("D7").Value = "'0.530"
$ws.Range("E7").Value = "  -2.00%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").Value = "'0.575"
$ws.Range("E9").Value = "  -2.88%  "

$ws.Range("D10").Value = "'35.75"
$ws.Range("E10").Value = "  -3.87%  "

$ws.Range("E11").Value = "  -0.95%  "

$ws.Range("D12").Value = "'0.0838"
$ws.Range("E12").Value = "  -0.65%  "

$ws.Range("D13").Value = "3.379.68"
$ws.Range("E13").Value = "  -1.76%  "

$ws.Range("D14").Value = "'17.81"
$ws.Range("E14").Value = "  -3.27%  "

$ws.Range("D15").Value = "'7.38"
$ws.Range("E15").Value = "  -2.49%  "

$ws.Range("B16").Value = "Uniswap"
$ws.Range("C16").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D16").Value = "'11.29"
$ws.Range("E16").Value = "  +52.73%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.915.00"
$ws.Range("E17").Value = "  -2.02%  "

$ws.Range("D18").Value = "'0.956"
$ws.Range("E18").Value = "  -1.06%  "

$ws.Range("D19").Value = "50.914.88"
$ws.Range("E19").Value = "  -1.03%  "

$ws.Range("D20").Value = "'3.11"
$ws.Range("E20").Value = "  -6.19%  "

$ws.Range("D21").Value = "'12.19"
$ws.Range("E21").Value = "  -5.60%  "

$ws.Range("D22").Value = "0.0₃0946"
$ws.Range("E22").Value = "  -1.61%  "

$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").Value = "'262.53"
$ws.Range("E23").Value = "  +0.31%  "

$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "'68.11"
$ws.Range("E24").Value = "  -1.90%  "

$ws.Range("D25").Value = "'3.09"
$ws.Range("E25").Value = "  +9.36%  "

$ws.Range("D26").Value = "'8.03"
$ws.Range("E26").Value = "  -2.60%  "

$ws.Range("D27").Value = "'7.18"
$ws.Range("E27").Value = "  -4.82%  "

$ws.Range("E28").Value = "  -0.02%  "

$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'25.35"
$ws.Range("E29").Value = "  -2.08%  "

$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").Value = "'0.161"
$ws.Range("E30").Value = "  -5.78%  "

$ws.Range("D31").Value = "'0.109"
$ws.Range("E31").Value = "  -3.08%  "

$ws.Range("D32").Value = "'9.81"
$ws.Range("E32").Value = "  -0.94%  "

$ws.Range("D33").Value = "'50.49"
$ws.Range("E33").Value = "  -0.99%  "

$ws.Range("E34").Value = "  -2.68%  "

$ws.Range("D35").Value = "'32.65"
$ws.Range("E35").Value = "  -6.77%  "

$ws.Range("D36").Value = "'0.0436"
$ws.Range("E36").Value = "  -2.14%  "

$ws.Range("E37").Value = "  -0.21%  "

$ws.Range("D38").Value = "'3.08"
$ws.Range("E38").Value = "  +1.45%  "

$ws.Range("E39").Value = "  -1.25%  "

$ws.Range("D40").Value = "'16.08"
$ws.Range("E40").Value = "  -6.26%  "

$ws.Range("E41").Value = "  -4.74%  "

$ws.Range("D42").Value = "'2.43"
$ws.Range("E42").Value = "  -5.92%  "

$ws.Range("D43").Value = "'119.24"
$ws.Range("E43").Value = "  -4.62%  "

$ws.Range("D44").Value = "'20.84"
$ws.Range("E44").Value = "  -3.87%  "

$ws.Range("E45").Value = "  -1.63%  "

$ws.Range("D46").Value = "'0.270"
$ws.Range("E46").Value = "  -6.96%  "

$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "'3.20"
$ws.Range("E47").Value = "  -0.49%  "

$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").Value = "'2.28"
$ws.Range("E48").Value = "  -3.47%  "

$ws.Range("D49").Value = "1.972.05"
$ws.Range("E49").Value = "  -3.09%  "

$ws.Range("D50").Value = "'0.0322"
$ws.Range("E50").Value = "  -5.45%  "

$ws.Range("D51").Value = "'1.29"
$ws.Range("E51").Value = "  +0.83%  "
